$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (shifts existing rows 11+ down by one, carrying styles/merges/heights)
$ws.Rows("11:11").Insert()

# Re-create the merges for the newly inserted row 11 (Insert() doesn't duplicate them automatically)
$ws.Range("B11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()

# Populate the newly inserted row with the new product entry
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "BONEDIVITON 50000L.U/2.5ML ORAL DROPS"
$ws.Cells.Item(11, 8).Value = "0:0"
$ws.Cells.Item(11, 12).Value = 80
$ws.Cells.Item(11, 14).Value = "1:0"

# Renumber the "م" sequence column (A) for every shifted product row (12 .. 58)
for ($r = 12; $r -le 58; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

# Update the running total (now on row 59) to include the new row's amount
$ws.Cells.Item(59, 11).Value = 3479.4499999999998
